$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @("Mark Thompson", "TRK-1015", "Downtown Construction", "06:00 AM", "03:00 PM", "N/A", "N/A", "On Time", "")
    3 = @("Christopher Thomas", "TRK-1011", "Downtown Construction", "06:30 AM", "03:30 PM", "06:29 AM", "03:33 PM", "On Time", "")
    4 = @("John Smith", "TRK-1001", "North Dallas Site", "06:30 AM", "03:30 PM", "06:25 AM", "03:35 PM", "On Time", "")
    5 = @("Richard Wilson", "TRK-1007", "Downtown Construction", "06:15 AM", "03:15 PM", "06:14 AM", "03:18 PM", "On Time", "")
    6 = @("William Brown", "TRK-1004", "Richardson Development", "06:30 AM", "03:30 PM", "06:32 AM", "02:54 PM", "Early End", "36 minutes early")
    7 = @("Thomas Taylor", "TRK-1009", "North Dallas Site", "07:00 AM", "04:00 PM", "06:43 AM", "03:52 PM", "On Time", "")
    8 = @("James Davis", "TRK-1005", "North Dallas Site", "06:00 AM", "03:00 PM", "07:03 AM", "03:12 PM", "Not On Job", "At incorrect location: North Richland Hills")
    9 = @("Joseph Moore", "TRK-1008", "Richardson Development", "06:30 AM", "03:30 PM", "06:55 AM", "03:25 PM", "Late", "25 minutes late")
    10 = @("Matthew Harris", "TRK-1013", "North Dallas Site", "06:45 AM", "03:45 PM", "06:52 AM", "03:37 PM", "On Time", "")
    11 = @("David Miller", "TRK-1006", "West Plano Project", "06:45 AM", "03:45 PM", "06:47 AM", "03:42 PM", "On Time", "")
    12 = @("Michael Johnson", "TRK-1002", "West Plano Project", "06:00 AM", "03:00 PM", "06:18 AM", "03:07 PM", "Late", "18 minutes late")
    13 = @("Daniel Jackson", "TRK-1012", "Richardson Development", "06:15 AM", "03:15 PM", "06:16 AM", "03:12 PM", "On Time", "")
    14 = @("Robert Williams", "TRK-1003", "Downtown Construction", "07:00 AM", "04:00 PM", "07:28 AM", "03:45 PM", "Late", "28 minutes late")
    15 = @("Anthony Martin", "TRK-1014", "West Plano Project", "07:00 AM", "04:00 PM", "07:16 AM", "03:58 PM", "Late", "16 minutes late")
    16 = @("Charles Anderson", "TRK-1010", "West Plano Project", "06:00 AM", "03:00 PM", "06:05 AM", "02:48 PM", "On Time", "")
}

foreach ($r in $data.Keys) {
    $row = $data[$r]
    for ($i = 0; $i -lt $row.Length; $i++) {
        $col = $i + 1
        $ws.Cells.Item($r, $col).Value = $row[$i]
    }
}

Write-Output "Done updating rows 2-16"
